$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17: averages / extremes of SW(S*)/SW(OPT) (col N) and SC(S*)/SC(OPT) (col Z)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format the new summary block: bold, 12pt, vertically centered
$fmt = $ws.Range("A14")
$fmt.Font.Bold = $true
$fmt.Font.Size = 12
$fmt.VerticalAlignment = -4108
$fmt.Copy()
$ws.Range("A14:B17").PasteSpecial(-4122)

# Match the selection left behind in the source workbook
[void]$ws.Range("A14:B17").Select()

# Page setup tweaks present in the edited workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "done"
